# Article 94 is going live.
#
# Row 7 shows a rolling window of the 3 most-recent "blog" tiles
# (C7 = newest .. I7 = oldest). Publishing article 94 slides the window
# forward by one slot:
#   I7: ser: 91 -> ser: 92
#   E7: ser: 92 -> ser: 93
#   C7: ser: 93 -> ser: 94   (brand new article)
# The shared string for "ser: 91" thus becomes unused (dropped on save)
# and a new shared string for "ser: 94" is created.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 92"
$ws.Range("E7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 93"
$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 94"

# Scroll the view and move the selection onto the slot that now holds
# the freshly-published article (I7), matching the saved workbook view.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I7").Select()
